$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("B583").Value = "19:52:35"
$ws1.Range("C583").Value = "19:59"
$ws1.Range("D583").Value = "17_ROMERO"
$ws1.Range("E583").Value = 7
$ws1.Range("F583").Value = "LP1912"
$ws1.Range("G583").Value = "30/12/2025"

$ws1.Range("B584").Value = "19:52:35"
$ws1.Range("C584").Value = "20:00"
$ws1.Range("D584").Value = "14_ABASTO"
$ws1.Range("E584").Value = 8
$ws1.Range("F584").Value = "LP1912"
$ws1.Range("G584").Value = "30/12/2025"

$ws1.Range("B585").Value = "19:52:35"
$ws1.Range("C585").Value = "20:00"
$ws1.Range("D585").Value = "16_SANTA ANA"
$ws1.Range("E585").Value = 8
$ws1.Range("F585").Value = "LP1912"
$ws1.Range("G585").Value = "30/12/2025"

$ws1.Range("B586").Value = "19:52:35"
$ws1.Range("C586").Value = "20:07"
$ws1.Range("D586").Value = "10_OLMOS"
$ws1.Range("E586").Value = 15
$ws1.Range("F586").Value = "LP1912"
$ws1.Range("G586").Value = "30/12/2025"

$ws1.Range("B587").Value = "19:52:35"
$ws1.Range("C587").Value = "20:09"
$ws1.Range("D587").Value = "15_ABASTO"
$ws1.Range("E587").Value = 17
$ws1.Range("F587").Value = "LP1912"
$ws1.Range("G587").Value = "30/12/2025"

$ws1.Range("B588").Value = "19:52:35"
$ws1.Range("C588").Value = "20:10"
$ws1.Range("D588").Value = "16_P MOR-167 Y 521"
$ws1.Range("E588").Value = 18
$ws1.Range("F588").Value = "LP1912"
$ws1.Range("G588").Value = "30/12/2025"

$ws1.Range("B589").Value = "19:52:35"
$ws1.Range("C589").Value = "20:19"
$ws1.Range("D589").Value = "23_HERNANDEZ"
$ws1.Range("E589").Value = 27
$ws1.Range("F589").Value = "LP1912"
$ws1.Range("G589").Value = "30/12/2025"

$ws1.Range("B590").Value = "19:52:35"
$ws1.Range("C590").Value = "20:20"
$ws1.Range("D590").Value = "26_HERNANDEZ"
$ws1.Range("E590").Value = 28
$ws1.Range("F590").Value = "LP1912"
$ws1.Range("G590").Value = "30/12/2025"

$ws1.Range("B591").Value = "19:52:35"
$ws1.Range("C591").Value = "20:22"
$ws1.Range("D591").Value = "11_ETCHEVERRY"
$ws1.Range("E591").Value = 30
$ws1.Range("F591").Value = "LP1912"
$ws1.Range("G591").Value = "30/12/2025"

$ws1.Range("B592").Value = "19:52:35"
$ws1.Range("C592").Value = "20:22"
$ws1.Range("D592").Value = "16_SANTA ANA"
$ws1.Range("E592").Value = 30
$ws1.Range("F592").Value = "LP1912"
$ws1.Range("G592").Value = "30/12/2025"

$ws1.Range("B593").Value = "19:52:35"
$ws1.Range("C593").Value = "20:23"
$ws1.Range("D593").Value = "215A_EL PATO"
$ws1.Range("E593").Value = 31
$ws1.Range("F593").Value = "LP1912"
$ws1.Range("G593").Value = "30/12/2025"

$ws1.Range("B594").Value = "19:52:35"
$ws1.Range("C594").Value = "20:52"
$ws1.Range("D594").Value = "15_ABASTO"
$ws1.Range("E594").Value = 60
$ws1.Range("F594").Value = "LP1912"
$ws1.Range("G594").Value = "30/12/2025"

$ws1.Range("B595").Value = "19:52:35"
$ws1.Range("C595").Value = "20:57"
$ws1.Range("D595").Value = "23_HERNANDEZ"
$ws1.Range("E595").Value = 65
$ws1.Range("F595").Value = "LP1912"
$ws1.Range("G595").Value = "30/12/2025"

$ws1.Range("B596").Value = "19:52:35"
$ws1.Range("C596").Value = "21:04"
$ws1.Range("D596").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("E596").Value = 72
$ws1.Range("F596").Value = "LP1912"
$ws1.Range("G596").Value = "30/12/2025"

$ws1.Range("B597").Value = "19:52:35"
$ws1.Range("C597").Value = "21:07"
$ws1.Range("D597").Value = "215B_EL PATO"
$ws1.Range("E597").Value = 75
$ws1.Range("F597").Value = "LP1912"
$ws1.Range("G597").Value = "30/12/2025"

$ws1.Range("B598").Value = "19:52:35"
$ws1.Range("C598").Value = "21:20"
$ws1.Range("D598").Value = "26_HERNANDEZ"
$ws1.Range("E598").Value = 88
$ws1.Range("F598").Value = "LP1912"
$ws1.Range("G598").Value = "30/12/2025"

$ws1.Range("B599").Value = "19:52:35"
$ws1.Range("C599").Value = "21:22"
$ws1.Range("D599").Value = "15_ABASTO"
$ws1.Range("E599").Value = 90
$ws1.Range("F599").Value = "LP1912"
$ws1.Range("G599").Value = "30/12/2025"

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 19:52:46"
$ws1.Range("A3").Value = "Total filas: 598"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("B41").Value = "30/12/2025"
$ws2.Range("C41").Value = "19:52:35"
$ws2.Range("D41").Value = "20:23"
$ws2.Range("E41").Value = "215A_EL PATO"
$ws2.Range("F41").Value = 31
$ws2.Range("G41").Value = "LP1912"

$ws2.Range("B42").Value = "30/12/2025"
$ws2.Range("C42").Value = "19:52:35"
$ws2.Range("D42").Value = "21:07"
$ws2.Range("E42").Value = "215B_EL PATO"
$ws2.Range("F42").Value = 75
$ws2.Range("G42").Value = "LP1912"

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 19:52:46"
$ws2.Range("A3").Value = "Total filas: 41"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("B74").Value = "30/12/2025"
$ws3.Range("C74").Value = "19:52:40"
$ws3.Range("D74").Value = "21:28"
$ws3.Range("E74").Value = "215C_LA PLATA"
$ws3.Range("F74").Value = 96
$ws3.Range("G74").Value = "L6203"

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 19:52:46"
$ws3.Range("A3").Value = "Total filas: 73"
